$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.Text = $tr.Text -replace [regex]::Escape("Huen Oh()"), "Huen Oh(301082798)"
